$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.021.77"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.038.79"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.94"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.84"
$ws.Range("E7").Value = "  +3.66%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0818"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("E11").Value = "  +0.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.67"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.339.10"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.45"
$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.767"
$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.19"
$ws.Range("E16").Value = "  -1.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.036.90"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.981.98"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.02"
$ws.Range("E19").Value = "  +0.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.98"
$ws.Range("E20").Value = "  -4.79%  "

$ws.Range("E21").Value = "  -0.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.98"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.38"
$ws.Range("E26").Value = "  +0.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.33"
$ws.Range("E27").Value = "  +0.45%  "

$ws.Range("E28").Value = "  -2.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.96"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.30"
$ws.Range("E30").Value = "  -2.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +1.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.16"
$ws.Range("E32").Value = "  +4.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.43"
$ws.Range("E33").Value = "  -1.85%  "

$ws.Range("E34").Value = "  -0.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.53"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.51"
$ws.Range("E36").Value = "  +6.78%  "

$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +0.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.527.13"
$ws.Range("E40").Value = "  +2.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.42"
$ws.Range("E41").Value = "  +5.58%  "

$ws.Range("E42").Value = "  +0.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.53"
$ws.Range("E43").Value = "  -0.37%  "

$ws.Range("E44").Value = "  -1.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0917"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").Value = "  -2.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.99"
$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.13"
$ws.Range("E50").Value = "  +0.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.227.78"
$ws.Range("E51").Value = "  -0.68%  "
